# Insert a new weekly record as row 14, shifting existing rows 14-32 down to 15-33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; this shifts rows 14:32 down to 15:33
# and copies formatting (incl. the date style in column D) from the row above.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly data entry.
$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 45225
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100107
$ws.Cells.Item(14, 8).Value = "Otros"
$ws.Cells.Item(14, 9).Value = 100107002
$ws.Cells.Item(14, 10).Value = "Chirimoya"
$ws.Cells.Item(14, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 80
$ws.Cells.Item(14, 14).Value = 21000
$ws.Cells.Item(14, 15).Value = 21000
$ws.Cells.Item(14, 16).Value = 21000
$ws.Cells.Item(14, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(14, 19).Value = 2100
$ws.Cells.Item(14, 20).Value = 10
